$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()
$excel.Goto($ws.Range("A5:B5"), $false) | Out-Null
Write-Output ("goto: " + $excel.Selection.Address() + " / " + $excel.ActiveCell.Address())
